$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: re-roll a volatile RANDBETWEEN formula cell until it lands on the
# exact target value, then deterministically recompute its dependents from
# the (now fixed) value so the SUM()/product formulas downstream don't get
# re-randomised by a fresh Calculate() pass.
# ---------------------------------------------------------------------------
function Set-RandCell($ws, $targetC5) {
    $found = $false
    $count = 0
    while (-not $found) {
        $ws.Range("C5").Formula = "=RANDBETWEEN(-100,100) / 100"
        $v = $ws.Range("C5").Value()
        $count = $count + 1
        if ($v -eq $targetC5) {
            $found = $true
        }
        if ($count -gt 200000) {
            break
        }
    }
    $ws.Range("D5").Formula = "=B5*C5"
    $ws.Range("G2").Formula = "=SUM(D2:D5)"
}

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)
$ws6 = $wb.Worksheets.Item(6)

# New "Airbus"/"Biomerieux" data: re-rolled RANDBETWEEN sensitivities on each
# day's sheet, cascading into the B5*C5 score and the SUM(D2:D5) total.
Set-RandCell $ws1 -0.87
Set-RandCell $ws2 -0.38
Set-RandCell $ws3 0.2
Set-RandCell $ws4 -0.92
Set-RandCell $ws5 -0.2
Set-RandCell $ws6 0.06

# Column A on the "27 02 2017" sheet is widened (best-fit) to its content.
$ws6.Columns.Item(1).AutoFit()
$ws6.Columns.Item(1).ColumnWidth = 20.59

# Selection / active-tab bookkeeping: sheet "27 02 2017" (6th tab) becomes
# the active one, scrolled back to the top, with Q3 selected; sheet
# "20 02 2017" keeps cell H30 selected but is no longer the active tab.
$ws1.Activate()
$ws1.Range("H30").Select()

$ws6.Activate()
$ws6.Range("Q3").Select()
